$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.270.46'
$ws.Range("E2").Value = '  -5.88%  '
$ws.Range("D3").Value = '1.668.70'
$ws.Range("E3").Value = '  -4.01%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.004'
$cell.ClearFormats()
$ws.Range("E4").Value = '  +0.30%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '217.61'
$cell.ClearFormats()
$ws.Range("E5").Value = '  -3.90%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.5078'
$cell.ClearFormats()
$ws.Range("E6").Value = '  -12.07%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '1.005'
$cell.ClearFormats()
$ws.Range("E7").Value = '  +0.24%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.2653'
$cell.ClearFormats()
$ws.Range("E8").Value = '  -3.10%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.06354'
$cell.ClearFormats()
$ws.Range("E9").Value = '  -4.17%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '21.49'
$cell.ClearFormats()
$ws.Range("E10").Value = '  -6.92%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07359'
$cell.ClearFormats()
$ws.Range("E11").Value = '  -2.53%  '
$ws.Range("D12").Value = '1.674.04'
$ws.Range("E12").Value = '  -3.92%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.543'
$cell.ClearFormats()
$ws.Range("E13").Value = '  -3.47%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.5808'
$cell.ClearFormats()
$ws.Range("E14").Value = '  -3.64%  '
$ws.Range("D15").Value = '1.894.12'
$ws.Range("E15").Value = '  -4.11%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.000008513'
$cell.ClearFormats()
$ws.Range("E16").Value = '  -2.07%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '64.72'
$cell.ClearFormats()
$ws.Range("E17").Value = '  -13.15%  '
$ws.Range("D18").Value = '26.323.47'
$ws.Range("E18").Value = '  -5.66%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '4.921'
$cell.ClearFormats()
$ws.Range("E19").Value = '  -7.36%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '1.005'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +0.35%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '10.83'
$cell.ClearFormats()
$ws.Range("E21").Value = '  -3.96%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '188.66'
$cell.ClearFormats()
$ws.Range("E22").Value = '  -7.98%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '6.187'
$cell.ClearFormats()
$ws.Range("E23").Value = '  -6.61%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '1.006'
$cell.ClearFormats()
$ws.Range("E24").Value = '  +0.30%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '143.26'
$cell.ClearFormats()
$ws.Range("E25").Value = '  -4.67%  '
$ws.Range("E26").Value = '  -4.52%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.1173'
$cell.ClearFormats()
$ws.Range("E27").Value = '  -4.91%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '15.63'
$cell.ClearFormats()
$ws.Range("E28").Value = '  -3.35%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.05847'
$cell.ClearFormats()
$ws.Range("E29").Value = '  -5.80%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '1.269'
$cell.ClearFormats()
$ws.Range("E30").Value = '  -8.23%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.324'
$cell.ClearFormats()
$ws.Range("E31").Value = '  -4.98%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.526'
$cell.ClearFormats()
$ws.Range("E32").Value = '  -5.68%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '3.510'
$cell.ClearFormats()
$ws.Range("E33").Value = '  -6.16%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.638'
$cell.ClearFormats()
$ws.Range("E34").Value = '  -2.61%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.011'
$cell.ClearFormats()
$ws.Range("E35").Value = '  -2.46%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.5992'
$cell.ClearFormats()
$ws.Range("E36").Value = '  -6.44%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '2.354'
$cell.ClearFormats()
$ws.Range("E37").Value = '  -3.66%  '
$ws.Range("E38").Value = '  -2.47%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.01612'
$cell.ClearFormats()
$ws.Range("E39").Value = '  -3.51%  '
$ws.Range("E40").Value = '  -2.53%  '
$ws.Range("D41").Value = '1.072.93'
$ws.Range("E41").Value = '  -4.24%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.8646'
$cell.ClearFormats()
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("E43").Value = '  +0.45%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '99.60'
$cell.ClearFormats()
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '1.815.31'
$ws.Range("E45").Value = '  -3.84%  '
$ws.Range("E46").Value = '  +0.88%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '55.73'
$cell.ClearFormats()
$ws.Range("E47").Value = '  -6.08%  '
$ws.Range("E48").Value = '  +0.19%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '8.083'
$cell.ClearFormats()
$ws.Range("E49").Value = '  -2.39%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.4292'
$cell.ClearFormats()
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("E51").Value = '  -3.68%  '

